$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Apendix3")
$ws4 = $wb.Worksheets.Item("Apendix4")

# Add the new "CMND" header column (F1) to both Apendix3 and Apendix4 sheets.
$ws3.Range("F1").Value = "CMND"
$ws3.Range("F1").Font.Bold = $true

$ws4.Range("F1").Value = "CMND"
$ws4.Range("F1").Font.Bold = $true
$ws4.PageSetup.Orientation = 1

# Update the selection on each sheet to the newly added cell, and make
# Apendix3 the active / selected tab (was Apendix4 before the edit).
[void]$ws4.Range("F1").Select()
[void]$ws3.Range("F1").Select()
[void]$ws3.Activate()
